$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B1 formula to reference $B2 instead of $C2 (will now evaluate to #VALUE! error)
$ws.Range("B1").Formula = '=TEXT(VALUE($B2),"ddd")'

# Add a new formula in G10 referencing an undefined name (produces #NAME? error)
$ws.Range("G10").Formula = "=av"

# Update B22 formula to reference $B21 instead of $B20
$ws.Range("B22").Formula = '=IFERROR($B21/$B5,"")'

# Update the selected/active cell for the frozen (bottom-right) pane
$ws.Range("B22").Select()
